# Insert a new data row at row 4 (pushes old rows 4..108 down to 5..109)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "Vega Monumental Concepción"
$ws.Range("C4").Value = "Bíobío"
$ws.Range("D4").Value = 45043
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 100112037
$ws.Range("G4").Value = "Cebollín"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 3300
$ws.Range("L4").Value = 3500
$ws.Range("M4").Value = 3389
$ws.Range("N4").Value = "`$/paquete 36 unidades"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 94
$ws.Range("Q4").Value = 36
$ws.Range("R4").Value = "Hortaliza"
